$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cost values for rows 2-60 (column L)
$newCosts = @(9,13,19,28,40,59,87,127,186,273,399,585,856,1254,1837,2691,3941,5772,8454,12382,18135,26561,38902,56977,83451,122226,179017,262195,384022,562454,823793,1206561,1767179,2588282,3790904,5552314,8132146,11910675,17444863,25550461,37422253,54810166,80277216,117577302,172208537,252223684,369417149,541063502,792463789,1160674956,1699972128,2000000000,3646734750,5341157231,7822877869,11457707664,16781428411,24578768089,35999071473)

for ($i = 0; $i -lt $newCosts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $newCosts[$i]
}

# Apply a numeric "0" format to the updated cost column range
$ws.Range("L2:L60").NumberFormat = "0"

# Update the active selection to match the edited column
$ws.Range("L2:L60").Select()
